# Update countries & provincias Spain
# Spain's numbers were updated and, as the table is sorted by total cases,
# Spain now ranks above Italy (rows 5/6 swap). Ethiopia's numbers were also
# updated and it now ranks above Bermudas/Guam (rows 146/147/148 shift).
# A couple of other rows only receive refreshed case counts, and the
# "last updated" timestamp string changes too.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 3 de Abril de 2020 a las 11:50"

# --- Estados Unidos (row 4) keeps its position; only stats refresh ---
$ws.Range("B4").Value = 245380
$ws.Range("C4").Value = 503
$ws.Range("D4").Value = 10411
$ws.Range("E4").Value = 228874
$ws.Range("F4").Value = 5421
$ws.Range("G4").Value = 25
$ws.Range("H4").Value = 6095

# --- España moves up to row 5 with refreshed stats ---
$ws.Range("A5").Value = "España"
$ws.Range("B5").Value = 117710
$ws.Range("C5").Value = 5645
$ws.Range("D5").Value = 30513
$ws.Range("E5").Value = 76262
$ws.Range("F5").Value = 6092
$ws.Range("G5").Value = 587
$ws.Range("H5").Value = 10935

# --- Italia moves down to row 6, keeping its previous stats ---
$ws.Range("A6").Value = "Italia"
$ws.Range("B6").Value = 115242
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 18278
$ws.Range("E6").Value = 83049
$ws.Range("F6").Value = 4053
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 13915

# --- Noruega (row 24) keeps its position; stats refresh ---
$ws.Range("B24").Value = 5255
$ws.Range("C24").Value = 108
$ws.Range("E24").Value = 5172
$ws.Range("G24").Value = 1
$ws.Range("H24").Value = 51

# --- Malasia (row 30): only "Muertes hoy" (F) refreshes ---
$ws.Range("F30").Value = 108

# --- Etiopia moves up to row 146 with refreshed stats ---
$ws.Range("A146").Value = "Etiopia"
$ws.Range("B146").Value = 35
$ws.Range("C146").Value = 6
$ws.Range("D146").Value = 3
$ws.Range("E146").Value = 32
$ws.Range("F146").Value = 2
$ws.Range("G146").Value = 0
$ws.Range("H146").Value = 0

# --- Bermudas moves down to row 147, keeping its previous stats ---
$ws.Range("A147").Value = "Bermudas"
$ws.Range("B147").Value = 35
$ws.Range("C147").Value = 0
$ws.Range("D147").Value = 11
$ws.Range("E147").Value = 24
$ws.Range("F147").Value = 0
$ws.Range("G147").Value = 0
$ws.Range("H147").Value = 0

# --- Guam moves down to row 148, keeping its previous stats ---
$ws.Range("A148").Value = "Guam"
$ws.Range("B148").Value = 32
$ws.Range("C148").Value = 0
$ws.Range("D148").Value = 0
$ws.Range("E148").Value = 31
$ws.Range("F148").Value = 0
$ws.Range("G148").Value = 0
$ws.Range("H148").Value = 1
